$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# A3: change value from "RO.ACT.001.CRE" to "RO.ACT.004EMP.SRL"
$ws.Range("A3").Value = "RO.ACT.004EMP.SRL"

# A4: add new value "AD.SEC.002.FON.01"
$ws.Range("A4").Value = "AD.SEC.002.FON.01"

# Update the selection to A10
$ws.Range("A10").Select()
